# Commit: "added details, delete and edit view"
#
# This adds new localization rows to the "Deutsch" resource sheet for the
# Pokemon Add/Details views as well as the experience-type enum, and
# renames two existing German translations (Basis Attacke -> Basis Angriff,
# Basis Spezialattacke -> Basis Spezialangriff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for the "Pokemon Add" view (rows 80-83)
$ws.Cells.Item(80,1).Value = "Pokemon_Add_Title"
$ws.Cells.Item(80,2).Value = "Pokémon hinzufügen"
$ws.Cells.Item(81,1).Value = "Pokemon_Add_Submit"
$ws.Cells.Item(81,2).Value = "Hinzufügen"
$ws.Cells.Item(82,1).Value = "Pokemon_Add_Cancel"
$ws.Cells.Item(82,2).Value = "Abbrechen"
$ws.Cells.Item(83,1).Value = "Pokemon_Add_BaseValues"
$ws.Cells.Item(83,2).Value = "Basiswerte"

# Rename existing base-attack / base-special-attack translations
$ws.Cells.Item(70,2).Value = "Basis Angriff"
$ws.Cells.Item(71,2).Value = "Basis Spezialangriff"

# New rows for "Pokemon Details" view and the Xp-type enum (rows 84-92)
$ws.Cells.Item(84,1).Value = "Pokemon_Details_Title"
$ws.Cells.Item(84,2).Value = "Details"
$ws.Cells.Item(85,1).Value = "Pokemon_Details_BaseValues"
$ws.Cells.Item(85,2).Value = "Basiswerte"
$ws.Cells.Item(86,1).Value = "Enums_XpType_Fast"
$ws.Cells.Item(86,2).Value = "Schnell"
$ws.Cells.Item(87,1).Value = "Enums_XpType_MediumFast"
$ws.Cells.Item(87,2).Value = "Mittel-Schnell"
$ws.Cells.Item(88,1).Value = "Enums_XpType_MediumSlow"
$ws.Cells.Item(88,2).Value = "Mittel-Langsam"
$ws.Cells.Item(89,1).Value = "Enums_XpType_Slow"
$ws.Cells.Item(89,2).Value = "Langsam"
$ws.Cells.Item(90,1).Value = "Enums_XpType_Erratic"
$ws.Cells.Item(90,2).Value = "Erratic"
$ws.Cells.Item(91,1).Value = "Enums_XpType_Fluctuating"
$ws.Cells.Item(91,2).Value = "Fluctuating"
$ws.Cells.Item(92,1).Value = "Pokemon_List_ExperienceType"
$ws.Cells.Item(92,2).Value = "Erfahrungstyp"

# Update the view's active cell / selection to match the end of the table
$ws.Cells.Item(93,2).Select()
